$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant
$xlPasteFormats = -4122

# --- Add new row 21 data -------------------------------------------------
# Order of the .Value assignments matters: it controls the order new
# strings are appended to the shared-strings table, which must match:
#   28 Renesas-RA6T2 Control Board (B21)
#   29 EPC9147D                    (C21)
#   30 RA6T2-EPC9194-DummyNema34-20k-2000n (F21)
#   31 EPC9194 Rev1_0              (D21)

$ws.Range("B21").Value = "Renesas-RA6T2 Control Board"
$ws.Range("B20").Copy()
$ws.Range("B21").PasteSpecial($xlPasteFormats)

$ws.Range("C21").Value = "EPC9147D"
$ws.Range("B20").Copy()
$ws.Range("C21").PasteSpecial($xlPasteFormats)

$ws.Range("F21").Value = "RA6T2-EPC9194-DummyNema34-20k-2000n"
$ws.Range("F9").Copy()
$ws.Range("F21").PasteSpecial($xlPasteFormats)

$ws.Range("D21").Value = "EPC9194 Rev1_0"
$ws.Range("F9").Copy()
$ws.Range("D21").PasteSpecial($xlPasteFormats)

$ws.Range("E21").Value = "Teknic_M-3411P-LN-08D"
$ws.Range("E20").Copy()
$ws.Range("E21").PasteSpecial($xlPasteFormats)

# --- Hyperlinks: F21 first (-> rId40), then D21 (-> rId41) ---------------
$ws.Hyperlinks.Add($ws.Range("F21"), "https://github.com/MarcoMacP/EPC-Reference-Designs-Firmware/tree/main/MOTOR%20DRIVE/RA6T2-EPC9194-DummyNema34-20k-2000n")
$ws.Range("F9").Copy()
$ws.Range("F21").PasteSpecial($xlPasteFormats)

$ws.Hyperlinks.Add($ws.Range("D21"), "https://epc-co.com/epc/products/demo-boards/epc9194")
$ws.Range("F9").Copy()
$ws.Range("D21").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# Leave the same cell selected as in the authored workbook
$ws.Range("E27").Select() | Out-Null
